$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert three new rows before current row 2 ("/public_facility") to make
# room for the new NULL/NULL row, a NULL row and a duplicated /public_facility row.
$ws.Rows("2:4").Insert()

# New row 2: NULL / NULL
$ws.Range("A2").Value = "NULL"
$ws.Range("B2").Value = "NULL"

# New row 3: /public_facility / NULL
$ws.Range("A3").Value = "/public_facility"
$ws.Range("B3").Value = "NULL"

# New row 4: NULL / 112399_sakado_city
$ws.Range("A4").Value = "NULL"
$ws.Range("B4").Value = "112399_sakado_city"

# Row 5 (previously row 2) already contains /public_facility / 112399_sakado_city.

$ws.Range("A9").Select()
